# Update coinranking.com crypto price/volume snapshot (GitHub Actions refresh).
# Values are written as literal text (matching the sheet's existing inline-string
# cells for columns D/E) rather than numbers, so figures like "14.00" or "0.0784"
# keep their exact printed form instead of being normalised to 14 / 0.0784 -> 14.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold text, not numbers. Any new D value that
# parses as a plain decimal (e.g. "292.45", "14.00") would otherwise be silently
# re-typed as a Number by Excel, dropping trailing zeros / losing the text format,
# so those specific cells are pre-formatted as Text before the value is written.
$textForcedCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D14", "D16", "D18", "D21", "D22", "D23", "D24", "D27", "D28", "D30", "D31", "D32", "D34", "D35", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '39.702.28'
$ws.Range("E2").Value = '  +1.22%  '
# Row 3
$ws.Range("D3").Value = '2.211.77'
$ws.Range("E3").Value = '  +0.95%  '
# Row 4
$ws.Range("E4").Value = '  +0.08%  '
# Row 5
$ws.Range("D5").Value = '292.45'
# Row 6
$ws.Range("D6").Value = '86.02'
$ws.Range("E6").Value = '  +5.80%  '
# Row 7
$ws.Range("E7").Value = '  +1.12%  '
# Row 8
$ws.Range("E8").Value = '  +0.01%  '
# Row 9
$ws.Range("D9").Value = '0.471'
$ws.Range("E9").Value = '  +1.49%  '
# Row 10
$ws.Range("D10").Value = '30.61'
$ws.Range("E10").Value = '  +5.58%  '
# Row 11
$ws.Range("D11").Value = '0.0784'
$ws.Range("E11").Value = '  +2.10%  '
# Row 12
$ws.Range("D12").Value = '47.50'
$ws.Range("E12").Value = '  +0.79%  '
# Row 13
$ws.Range("E13").Value = '  +1.47%  '
# Row 14
$ws.Range("D14").Value = '6.33'
$ws.Range("E14").Value = '  +1.86%  '
# Row 15
$ws.Range("D15").Value = '2.554.56'
$ws.Range("E15").Value = '  +1.06%  '
# Row 16
$ws.Range("D16").Value = '14.00'
$ws.Range("E16").Value = '  +0.70%  '
# Row 17
$ws.Range("D17").Value = '2.219.13'
$ws.Range("E17").Value = '  +1.46%  '
# Row 18
$ws.Range("D18").Value = '0.729'
$ws.Range("E18").Value = '  +3.14%  '
# Row 19
$ws.Range("D19").Value = '39.680.70'
$ws.Range("E19").Value = '  +1.47%  '
# Row 20
$ws.Range("D20").Value = '0.0₃0879'
$ws.Range("E20").Value = '  +1.44%  '
# Row 21
$ws.Range("D21").Value = '11.22'
$ws.Range("E21").Value = '  +10.03%  '
# Row 22
$ws.Range("D22").Value = '5.80'
$ws.Range("E22").Value = '  +2.00%  '
# Row 23
$ws.Range("D23").Value = '65.45'
$ws.Range("E23").Value = '  +1.08%  '
# Row 24
$ws.Range("D24").Value = '235.19'
$ws.Range("E24").Value = '  +4.62%  '
# Row 26
$ws.Range("E26").Value = '  +2.81%  '
# Row 27
$ws.Range("D27").Value = '1.83'
$ws.Range("E27").Value = '  +2.36%  '
# Row 28
$ws.Range("D28").Value = '22.71'
$ws.Range("E28").Value = '  +1.23%  '
# Row 29
$ws.Range("E29").Value = '  +1.54%  '
# Row 30
$ws.Range("D30").Value = '9.22'
$ws.Range("E30").Value = '  +2.10%  '
# Row 31
$ws.Range("D31").Value = '32.77'
$ws.Range("E31").Value = '  +3.84%  '
# Row 32
$ws.Range("D32").Value = '151.51'
$ws.Range("E32").Value = '  +1.37%  '
# Row 33
$ws.Range("E33").Value = '  -0.12%  '
# Row 34
$ws.Range("D34").Value = '4.92'
$ws.Range("E34").Value = '  +2.79%  '
# Row 35
$ws.Range("D35").Value = '0.0717'
$ws.Range("E35").Value = '  +3.87%  '
# Row 36
$ws.Range("E36").Value = '  +2.03%  '
# Row 37
$ws.Range("E37").Value = '  +6.96%  '
# Row 38
$ws.Range("D38").Value = '0.111'
$ws.Range("E38").Value = '  +1.59%  '
# Row 39
$ws.Range("D39").Value = '15.84'
$ws.Range("E39").Value = '  +3.67%  '
# Row 40
$ws.Range("D40").Value = '0.0986'
$ws.Range("E40").Value = '  +2.84%  '
# Row 41
$ws.Range("D41").Value = '1.69'
$ws.Range("E41").Value = '  +4.14%  '
# Row 42
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.067.03'
$ws.Range("E42").Value = '  +9.43%  '
# Row 43
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '3.77'
$ws.Range("E43").Value = '  +5.29%  '
# Row 44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0267'
$ws.Range("E44").Value = '  +3.31%  '
# Row 45
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '9.94'
$ws.Range("E45").Value = '  +11.58%  '
# Row 46
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '2.09'
$ws.Range("E46").Value = '  +0.34%  '
# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '17.65'
$ws.Range("E47").Value = '  +10.05%  '
# Row 48
$ws.Range("D48").Value = '2.59'
$ws.Range("E48").Value = '  +0.11%  '
# Row 49
$ws.Range("D49").Value = '2.432.52'
$ws.Range("E49").Value = '  +1.43%  '
# Row 50
$ws.Range("D50").Value = '70.98'
$ws.Range("E50").Value = '  +0.10%  '
# Row 51
$ws.Range("D51").Value = '88.88'
$ws.Range("E51").Value = '  +2.57%  '
